# Robust test and unit time/cost assessment
# - Switch active sheet/tab to "constraints"
# - Update the "cost" column (G) with the revised unit cost figures
# - Move the sheet selection to G5 (first cost cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constraints")

# Revised unit cost values (column G, rows 5-16)
$ws.Range("G5").Value = 15000
$ws.Range("G6").Value = 78000
$ws.Range("G7").Value = 5000
$ws.Range("G8").Value = 120000
$ws.Range("G9").Value = 54000
$ws.Range("G10").Value = 100000
$ws.Range("G11").Value = 150000
$ws.Range("G12").Value = 65000
$ws.Range("G13").Value = 45000
$ws.Range("G14").Value = 120000
$ws.Range("G15").Value = 150000
$ws.Range("G16").Value = 43200

# Make "constraints" the active sheet/tab and select the first cost cell
$ws.Activate() | Out-Null
$ws.Range("G5").Select() | Out-Null
